# Fruta / hortaliza, semanal
# Insert a new daily price record for "Arándano (blue)" at row 64 of the
# consolidated sheet, pushing the existing records (old rows 64-158) down
# by one row (to 65-159).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 64 (shifts rows 64..158 down to 65..159).
$ws.Rows.Item(64).Insert()

# Populate the newly inserted row 64 with the new observation.
$ws.Cells.Item(64, 1).Value  = 9
$ws.Cells.Item(64, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(64, 3).Value  = "Metropolitana"
$ws.Cells.Item(64, 4).Value  = 44579
$ws.Cells.Item(64, 5).Value  = 13
$ws.Cells.Item(64, 6).Value  = "Fruta"
$ws.Cells.Item(64, 7).Value  = 100101
$ws.Cells.Item(64, 8).Value  = "Berries"
$ws.Cells.Item(64, 9).Value  = 100101001
$ws.Cells.Item(64, 10).Value = "Arándano (blue)"
$ws.Cells.Item(64, 11).Value = "Sin especificar"
$ws.Cells.Item(64, 12).Value = "Primera"
$ws.Cells.Item(64, 13).Value = 480
$ws.Cells.Item(64, 14).Value = 4000
$ws.Cells.Item(64, 15).Value = 4000
$ws.Cells.Item(64, 16).Value = 4000
$ws.Cells.Item(64, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(64, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(64, 19).Value = 2000
$ws.Cells.Item(64, 20).Value = 2
